$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct "ID" column (G) typo: values were stored as plain numbers (1-7),
# but they should be text ids formatted as "id01".."id07".
$ids = @("id01", "id02", "id03", "id04", "id05", "id06", "id07")

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 7).Value = $ids[$r - 2]
}
for ($r = 8; $r -le 14; $r++) {
    $ws.Cells.Item($r, 7).Value = $ids[$r - 8]
}

# Reflect the edited range as the active selection, like Excel would leave it
# after the user finished typing/pasting over G2:G14.
$ws.Range("G2:G14").Select()
